$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 52, shifting existing rows 52:85 down to 53:86
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new data record
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = 44767
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100108
$ws.Range("H52").Value = "Tropicales y subtropicales"
$ws.Range("I52").Value = 100108007
$ws.Range("J52").Value = "Coco"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 15
$ws.Range("N52").Value = 30000
$ws.Range("O52").Value = 30000
$ws.Range("P52").Value = 30000
$ws.Range("Q52").Value = '$/malla 20 unidades'
$ws.Range("R52").Value = "Perú"
$ws.Range("S52").Value = 1500
$ws.Range("T52").Value = 20
